# Applies the "comandosGit" update:
#   - adds a "Ver el registro..." / "git log" follow-up: a blank paragraph
#     and a new "mas.." paragraph at the end of the document
#   - the _GoBack bookmark (previously wrapping the end of the "git log"
#     paragraph) moves so it wraps the new "mas.." text instead

$d = $word.ActiveDocument

# The existing document ends with a paragraph containing a tab + "git log",
# and carries the (hidden) _GoBack bookmark right at the end of that
# paragraph. Remove it from there first so it can be re-added around the
# new trailing text.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Locate the "git log" paragraph (the last paragraph in the document).
$gitLogPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Insert a new, empty paragraph right after it.
$r = $gitLogPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Insert a further new paragraph after the blank one, and give it the
# "mas.." text.
$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $blankPara.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$masPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$masPara.Range.Text = "mas.."

# Re-create the _GoBack bookmark around the new last paragraph's text.
$d.Bookmarks.Add("_GoBack", $masPara.Range)
